{"js": "// Update the \"Description : \" intro paragraph with the full description text,\n// and move each data row's \"Commentaire\" text (last column) into the\n// \"Description\" column (5th column), leaving the \"Commentaire\" column empty.\n\nconst body = context.document.body;\n\n// --- 1. Update the \"Description : \" paragraph ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst descParagraph = paragraphs.items.find((p) => p.text === \"Description : \");\nif (descParagraph) {\n  descParagraph.clear();\n  descParagraph.insertText(\n    \"Description : Codes inter-sant\u00e9 uniquement. D\u00e9crit le cadre conventionnel dans lequel s'inscrit une demande de ressources, voir EMSI pour lien 15-Nexsis.\",\n    \"Start\"\n  );\n  await context.sync();\n}\n\n// --- 2. Move \"Commentaire\" text into \"Description\" for each data row ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The second table in the document is the Code/Libell\u00e9/.../Description/Commentaire table.\nconst tbl = tables.items[1];\ntbl.load(\"rowCount,values\");\nawait context.sync();\n\nconst DESCRIPTION_COL = 4; // 0-based index of \"Description\" column\nconst COMMENTAIRE_COL = 5; // 0-based index of \"Commentaire\" column\n\nfor (let row = 1; row < tbl.rowCount; row++) {\n  const commentaireText = tbl.values[row][COMMENTAIRE_COL];\n  if (!commentaireText) {\n    continue;\n  }\n\n  const descCell = tbl.getCell(row, DESCRIPTION_COL);\n  descCell.body.insertText(commentaireText, \"Replace\");\n\n  const commentCell = tbl.getCell(row, COMMENTAIRE_COL);\n  commentCell.body.clear();\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Description : \" intro paragraph with the full description text,\n# and move each data row's \"Commentaire\" text (last column) into the\n# \"Description\" column (5th column), leaving the \"Commentaire\" column empty.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the \"Description : \" paragraph ---\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  $p = $paras.Item($i)\n  if ($p.Range.Text -eq \"Description : `r\") {\n    $r = $p.Range\n    # Exclude the trailing paragraph mark so we don't spawn a new paragraph.\n    $subRng = $d.Range($r.Start, $r.End - 1)\n    $subRng.Text = \"Description : Codes inter-sant\u00e9 uniquement. D\u00e9crit le cadre conventionnel dans lequel s'inscrit une demande de ressources, voir EMSI pour lien 15-Nexsis.\"\n    break\n  }\n}\n\n# --- 2. Move \"Commentaire\" text into \"Description\" for each data row ---\n# The second table in the document is the Code/Libell\u00e9/.../Description/Commentaire table.\n$tbl = $d.Tables.Item(2)\n\n$descCol = 5\n$commentCol = 6\n\nfor ($row = 2; $row -le $tbl.Rows.Count; $row++) {\n  $commentCell = $tbl.Cell($row, $commentCol)\n  $commentText = $commentCell.Range.Text\n  # Strip the trailing cell-mark characters (carriage return + cell mark).\n  $commentText = $commentText.TrimEnd([char]7).TrimEnd([char]13)\n\n  if ([string]::IsNullOrEmpty($commentText)) {\n    continue\n  }\n\n  $descCell = $tbl.Cell($row, $descCol)\n  $descCell.Range.Text = $commentText\n  $commentCell.Range.Text = \"\"\n}\n"}
